$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the X/Y coordinate values in columns B and C (rows 2-42) ---
# (recomputed coordinates, e.g. after reprojecting the demand points)
$ws.Cells.Item(2, 2).Value = 30.5711296141195
$ws.Cells.Item(2, 3).Value = 37.7692791490538
$ws.Cells.Item(3, 2).Value = 30.5361826961962
$ws.Cells.Item(3, 3).Value = 37.7717753823126
$ws.Cells.Item(4, 2).Value = 30.5389666077598
$ws.Cells.Item(4, 3).Value = 37.7611887466813
$ws.Cells.Item(5, 2).Value = 30.5236449889832
$ws.Cells.Item(5, 3).Value = 37.7610028743893
$ws.Cells.Item(6, 2).Value = 30.529855974324
$ws.Cells.Item(6, 3).Value = 37.7798746832696
$ws.Cells.Item(7, 2).Value = 30.5575996303048
$ws.Cells.Item(7, 3).Value = 37.7704178033416
$ws.Cells.Item(8, 2).Value = 30.5844729616467
$ws.Cells.Item(8, 3).Value = 37.7540376592585
$ws.Cells.Item(9, 2).Value = 30.5447728310508
$ws.Cells.Item(9, 3).Value = 37.7678001831541
$ws.Cells.Item(10, 2).Value = 30.5472444521171
$ws.Cells.Item(10, 3).Value = 37.7790170893297
$ws.Cells.Item(11, 2).Value = 30.5691474941673
$ws.Cells.Item(11, 3).Value = 37.7622131223577
$ws.Cells.Item(12, 2).Value = 30.5141888528279
$ws.Cells.Item(12, 3).Value = 37.7716603420279
$ws.Cells.Item(13, 2).Value = 30.5774543601264
$ws.Cells.Item(13, 3).Value = 37.7512200401011
$ws.Cells.Item(14, 2).Value = 30.5688449427822
$ws.Cells.Item(14, 3).Value = 37.7565214036057
$ws.Cells.Item(15, 2).Value = 30.5511854390431
$ws.Cells.Item(15, 3).Value = 37.7662886866951
$ws.Cells.Item(16, 2).Value = 30.5567133701435
$ws.Cells.Item(16, 3).Value = 37.7527620179125
$ws.Cells.Item(17, 2).Value = 30.530804076422
$ws.Cells.Item(17, 3).Value = 37.7747011143654
$ws.Cells.Item(18, 2).Value = 30.5396320397447
$ws.Cells.Item(18, 3).Value = 37.7664832631085
$ws.Cells.Item(19, 2).Value = 30.5637164423968
$ws.Cells.Item(19, 3).Value = 37.7537541896793
$ws.Cells.Item(20, 2).Value = 30.5596908474654
$ws.Cells.Item(20, 3).Value = 37.7672806521465
$ws.Cells.Item(21, 2).Value = 30.5118347246543
$ws.Cells.Item(21, 3).Value = 37.8229782859424
$ws.Cells.Item(22, 2).Value = 30.5113974779971
$ws.Cells.Item(22, 3).Value = 37.7647096216254
$ws.Cells.Item(23, 2).Value = 30.5683081539031
$ws.Cells.Item(23, 3).Value = 37.7533517532827
$ws.Cells.Item(24, 2).Value = 30.5495385824677
$ws.Cells.Item(24, 3).Value = 37.7630222496169
$ws.Cells.Item(25, 2).Value = 30.5409785867431
$ws.Cells.Item(25, 3).Value = 37.8180415465428
$ws.Cells.Item(26, 2).Value = 30.5534791459276
$ws.Cells.Item(26, 3).Value = 37.7841949692717
$ws.Cells.Item(27, 2).Value = 30.5554681971292
$ws.Cells.Item(27, 3).Value = 37.7626308576564
$ws.Cells.Item(28, 2).Value = 30.5307863066809
$ws.Cells.Item(28, 3).Value = 37.7527476162981
$ws.Cells.Item(29, 2).Value = 30.546587789673
$ws.Cells.Item(29, 3).Value = 37.7560592687151
$ws.Cells.Item(30, 2).Value = 30.5357056027002
$ws.Cells.Item(30, 3).Value = 37.7868239178745
$ws.Cells.Item(31, 2).Value = 30.5740318801062
$ws.Cells.Item(31, 3).Value = 37.7584863666919
$ws.Cells.Item(32, 2).Value = 30.560003340186
$ws.Cells.Item(32, 3).Value = 37.7613615569701
$ws.Cells.Item(33, 2).Value = 30.5553639935451
$ws.Cells.Item(33, 3).Value = 37.7595079328201
$ws.Cells.Item(34, 2).Value = 30.5547456511584
$ws.Cells.Item(34, 3).Value = 37.7644310623444
$ws.Cells.Item(35, 2).Value = 30.5489605661489
$ws.Cells.Item(35, 3).Value = 37.7885773929111
$ws.Cells.Item(36, 2).Value = 30.5608980484508
$ws.Cells.Item(36, 3).Value = 37.777540282213
$ws.Cells.Item(37, 2).Value = 30.5640707909029
$ws.Cells.Item(37, 3).Value = 37.7631876423035
$ws.Cells.Item(38, 2).Value = 30.5796235696272
$ws.Cells.Item(38, 3).Value = 37.8209857603662
$ws.Cells.Item(39, 2).Value = 30.5563440888914
$ws.Cells.Item(39, 3).Value = 37.756448928897
$ws.Cells.Item(40, 2).Value = 30.5526467138055
$ws.Cells.Item(40, 3).Value = 37.7580589495252
$ws.Cells.Item(41, 2).Value = 30.5892163509704
$ws.Cells.Item(41, 3).Value = 37.7581912142887
$ws.Cells.Item(42, 2).Value = 30.5660173081617
$ws.Cells.Item(42, 3).Value = 37.7880536505561

# --- Widen columns B (X) and C (Y) so the longer decimal values are fully visible ---
$ws.Columns.Item(2).ColumnWidth = 18.2501
$ws.Columns.Item(3).ColumnWidth = 19.5834

# --- Move the selection down to C45 (scrolled towards the bottom of the list) ---
$ws.Range("C45").Select()
